# Auto-generated Excel COM-interop script to apply the Ultima_Profits market-data refresh.
# Updates currentAveragePrice*/LevePrice*/LeveProfit* columns (H:N) for the affected Leve rows
# on each crafting-class sheet, matching the latest Universalis market snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1666.9656
$ws.Range("I40").Value = 1600
$ws.Range("J40").Value = 1738.7142
$ws.Range("K40").Value = 1600
$ws.Range("L40").Value = 1738.7142
$ws.Range("M40").Value = -1425
$ws.Range("N40").Value = -2088.7142
$ws.Range("H108").Value = 54795
$ws.Range("J108").Value = 54795
$ws.Range("L108").Value = 54795
$ws.Range("N108").Value = -62475
$ws.Range("H112").Value = 1493.4482
$ws.Range("J112").Value = 1611.9231
$ws.Range("L112").Value = 4835.7693
$ws.Range("N112").Value = -7051.7693
$ws.Range("H127").Value = 726
$ws.Range("J127").Value = 875.5833
$ws.Range("L127").Value = 2626.7499
$ws.Range("N127").Value = -12546.7499
$ws.Range("H137").Value = 5264061
$ws.Range("I137").Value = 875.8182
$ws.Range("J137").Value = 12500941
$ws.Range("K137").Value = 2627.4546
$ws.Range("L137").Value = 37502823
$ws.Range("M137").Value = -77.45460000000003
$ws.Range("N137").Value = -37507923
$ws.Range("H138").Value = 2094.5146
$ws.Range("J138").Value = 3879.96
$ws.Range("L138").Value = 11639.88
$ws.Range("N138").Value = -21919.88

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 333.81818
$ws.Range("I4").Value = 423.33334
$ws.Range("J4").Value = 300.25
$ws.Range("K4").Value = 423.33334
$ws.Range("L4").Value = 300.25
$ws.Range("M4").Value = -307.33334
$ws.Range("N4").Value = -532.25
$ws.Range("H43").Value = 6803
$ws.Range("J43").Value = 6383.0835
$ws.Range("L43").Value = 6383.0835
$ws.Range("N43").Value = -7009.0835
$ws.Range("H61").Value = 12822052
$ws.Range("I61").Value = 16668115
$ws.Range("J61").Value = 1842.1111
$ws.Range("K61").Value = 16668115
$ws.Range("L61").Value = 1842.1111
$ws.Range("M61").Value = -16667903
$ws.Range("N61").Value = -2266.1111
$ws.Range("H74").Value = 11907044
$ws.Range("I74").Value = 16130226
$ws.Range("K74").Value = 16130226
$ws.Range("M74").Value = -16129352
$ws.Range("H77").Value = 11907044
$ws.Range("I77").Value = 16130226
$ws.Range("K77").Value = 80651130
$ws.Range("M77").Value = -80646762
$ws.Range("H109").Value = 32335.111
$ws.Range("J109").Value = 32335.111
$ws.Range("L109").Value = 32335.111
$ws.Range("N109").Value = -35109.111
$ws.Range("H110").Value = 1184.52
$ws.Range("I110").Value = 937.5263
$ws.Range("J110").Value = 1966.6666
$ws.Range("K110").Value = 937.5263
$ws.Range("L110").Value = 1966.6666
$ws.Range("M110").Value = 1107.4737
$ws.Range("N110").Value = -6056.6666
$ws.Range("H132").Value = 5815857.5
$ws.Range("I132").Value = 7814132
$ws.Range("K132").Value = 23442396
$ws.Range("M132").Value = -23439866
$ws.Range("H136").Value = 12822052
$ws.Range("I136").Value = 16668115
$ws.Range("J136").Value = 1842.1111
$ws.Range("K136").Value = 50004345
$ws.Range("L136").Value = 5526.3333
$ws.Range("M136").Value = -50001795
$ws.Range("N136").Value = -10626.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1027.6923
$ws.Range("I99").Value = 1028.1818
$ws.Range("K99").Value = 1028.1818
$ws.Range("M99").Value = 469.8181999999999
$ws.Range("H105").Value = 4995.2383
$ws.Range("J105").Value = 5117.647
$ws.Range("L105").Value = 5117.647
$ws.Range("N105").Value = -8611.647000000001
$ws.Range("H107").Value = 1247.6538
$ws.Range("I107").Value = 1144.1578
$ws.Range("J107").Value = 1528.5714
$ws.Range("K107").Value = 1144.1578
$ws.Range("L107").Value = 1528.5714
$ws.Range("M107").Value = 775.8422
$ws.Range("N107").Value = -5368.5714
$ws.Range("H134").Value = 2776.3462
$ws.Range("I134").Value = 1694.7727
$ws.Range("J134").Value = 8725
$ws.Range("K134").Value = 5084.3181
$ws.Range("L134").Value = 26175
$ws.Range("M134").Value = -2549.3181
$ws.Range("N134").Value = -31245

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5851465.5
$ws.Range("I31").Value = 4274.5527
$ws.Range("J31").Value = 17545848
$ws.Range("K31").Value = 4274.5527
$ws.Range("L31").Value = 17545848
$ws.Range("M31").Value = -3979.5527
$ws.Range("N31").Value = -17546438
$ws.Range("H34").Value = 5851465.5
$ws.Range("I34").Value = 4274.5527
$ws.Range("J34").Value = 17545848
$ws.Range("K34").Value = 4274.5527
$ws.Range("L34").Value = 17545848
$ws.Range("M34").Value = -4072.5527
$ws.Range("N34").Value = -17546252
$ws.Range("H132").Value = 10871149
$ws.Range("I132").Value = 12501124
$ws.Range("J132").Value = 4652
$ws.Range("K132").Value = 37503372
$ws.Range("L132").Value = 13956
$ws.Range("M132").Value = -37500842
$ws.Range("N132").Value = -19016

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1176
$ws.Range("I5").Value = 664.8
$ws.Range("J5").Value = 1310.5264
$ws.Range("K5").Value = 1994.4
$ws.Range("L5").Value = 3931.5792
$ws.Range("M5").Value = -1882.4
$ws.Range("N5").Value = -4155.5792
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H122").Value = 1301.3636
$ws.Range("J122").Value = 1999
$ws.Range("L122").Value = 17991
$ws.Range("N122").Value = -22891
$ws.Range("H132").Value = 1389
$ws.Range("I132").Value = 860
$ws.Range("J132").Value = 3505
$ws.Range("K132").Value = 7740
$ws.Range("L132").Value = 31545
$ws.Range("M132").Value = -5210
$ws.Range("N132").Value = -36605
$ws.Range("H135").Value = 1176
$ws.Range("I135").Value = 664.8
$ws.Range("J135").Value = 1310.5264
$ws.Range("K135").Value = 5983.2
$ws.Range("L135").Value = 11794.7376
$ws.Range("M135").Value = -3448.2
$ws.Range("N135").Value = -16864.7376

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5557840
$ws.Range("I122").Value = 7408787
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 22226361
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -22223911
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 3724.8096
$ws.Range("I132").Value = 2685.3125
$ws.Range("J132").Value = 7051.2
$ws.Range("K132").Value = 8055.9375
$ws.Range("L132").Value = 21153.6
$ws.Range("M132").Value = -5525.9375
$ws.Range("N132").Value = -26213.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1213
$ws.Range("I22").Value = 528.5714
$ws.Range("J22").Value = 1555.2142
$ws.Range("K22").Value = 528.5714
$ws.Range("L22").Value = 1555.2142
$ws.Range("M22").Value = -233.5714
$ws.Range("N22").Value = -2145.2142
$ws.Range("H27").Value = 1213
$ws.Range("I27").Value = 528.5714
$ws.Range("J27").Value = 1555.2142
$ws.Range("K27").Value = 528.5714
$ws.Range("L27").Value = 1555.2142
$ws.Range("M27").Value = -421.5714
$ws.Range("N27").Value = -1769.2142
$ws.Range("H46").Value = 990.9167
$ws.Range("I46").Value = 878.2
$ws.Range("K46").Value = 878.2
$ws.Range("M46").Value = -690.2
$ws.Range("H132").Value = 13524.883
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 13524.883
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 40574.649
$ws.Range("N132").Value = -45634.649
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 100000
$ws.Range("I74").Value = 100000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 100000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -99064
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 100000
$ws.Range("I77").Value = 100000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 300000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -295320
$ws.Range("N77").ClearContents()
$ws.Range("H122").Value = 1587.1852
$ws.Range("J122").Value = 866.125
$ws.Range("L122").Value = 2598.375
$ws.Range("N122").Value = -7498.375
$ws.Range("H132").Value = 4667.1665
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4667.1665
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 14001.4995
$ws.Range("N132").Value = -19061.4995
$ws.Range("M132").ClearContents()
